# Updated cryptos list values (Price column D, Volume(1h) column E)
# per the Tue Jun 13 19:30:04 UTC 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.885.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.740.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5180"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.62%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2745"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06149"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.741.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07167"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.6446"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.00%  "
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.600"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9996"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.912.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("E19").Value = "  +2.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006778"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.963.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.278"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.676"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.249"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.508"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.766"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.938"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08304"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.665"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04591"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.655"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9902"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6185"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.688"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.98%  "
$ws.Range("E38").Value = "  +3.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.933"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9991"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3841"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7387"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.980"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1127"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.214"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05264"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.28%  "
$ws.Range("E49").Value = "  +2.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.638"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("E51").Value = "  +1.25%  "
